$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$arrBC = New-Object 'object[,]' 24,2
$arrBC[0,0] = 11.96991345826797
$arrBC[0,1] = 10.81120567706507
$arrBC[1,0] = 11.6551591358343
$arrBC[1,1] = 10.83439586924669
$arrBC[2,0] = 11.45852491512911
$arrBC[2,1] = 10.84951521950415
$arrBC[3,0] = 11.37765811286709
$arrBC[3,1] = 10.85589844696599
$arrBC[4,0] = 11.36418913123472
$arrBC[4,1] = 10.85697180014103
$arrBC[5,0] = 11.4574371480591
$arrBC[5,1] = 10.84960040643958
$arrBC[6,0] = 11.86215132426379
$arrBC[6,1] = 10.8190192520527
$arrBC[7,0] = 12.62446302189016
$arrBC[7,1] = 10.76601044537274
$arrBC[8,0] = 13.15982433335573
$arrBC[8,1] = 10.7312732227158
$arrBC[9,0] = 13.39696831393677
$arrBC[9,1] = 10.71637677894351
$arrBC[10,0] = 13.48577208868957
$arrBC[10,1] = 10.71086556164922
$arrBC[11,0] = 13.46669213751891
$arrBC[11,1] = 10.71204673837325
$arrBC[12,0] = 13.40429466686242
$arrBC[12,1] = 10.71592077035474
$arrBC[13,0] = 13.36594231676199
$arrBC[13,1] = 10.71831060855497
$arrBC[14,0] = 13.14419091364676
$arrBC[14,1] = 10.73226492275199
$arrBC[15,0] = 13.00645806378694
$arrBC[15,1] = 10.74105705800425
$arrBC[16,0] = 12.92664008640697
$arrBC[16,1] = 10.74619933894366
$arrBC[17,0] = 12.89951493041309
$arrBC[17,1] = 10.74795508855738
$arrBC[18,0] = 13.02118240629754
$arrBC[18,1] = 10.74011229792654
$arrBC[19,0] = 13.42264997482028
$arrBC[19,1] = 10.71477935598623
$arrBC[20,0] = 13.67918679143603
$arrBC[20,1] = 10.69897882859657
$arrBC[21,0] = 13.54282706673061
$arrBC[21,1] = 10.70734285382849
$arrBC[22,0] = 13.01452749834802
$arrBC[22,1] = 10.74053915116472
$arrBC[23,0] = 12.4221991454838
$arrBC[23,1] = 10.77960919090583
$ws.Range("B2:C25").Value = $arrBC

$arrEL = New-Object 'object[,]' 24,8
$arrEL[0,0] = 16.57964772017429
$arrEL[0,1] = 36.77951474306486
$arrEL[0,2] = 30.13312868227559
$arrEL[0,3] = 14.81451105515125
$arrEL[0,4] = 23.88884783221761
$arrEL[0,5] = 7.763676358809629
$arrEL[0,6] = 8.316685527043457
$arrEL[0,7] = 12.70137982175408
$arrEL[1,0] = 16.57563495243377
$arrEL[1,1] = 36.8313907054996
$arrEL[1,2] = 30.29311796827507
$arrEL[1,3] = 14.86851405924421
$arrEL[1,4] = 23.99752384844511
$arrEL[1,5] = 7.751039829416176
$arrEL[1,6] = 8.087301579925686
$arrEL[1,7] = 12.67082374872037
$arrEL[2,0] = 16.57574010105909
$arrEL[2,1] = 36.87156926367057
$arrEL[2,2] = 30.3997184747962
$arrEL[2,3] = 14.90372665277077
$arrEL[2,4] = 24.06834997366074
$arrEL[2,5] = 7.743226244150094
$arrEL[2,6] = 7.943472941502468
$arrEL[2,7] = 12.65374400608392
$arrEL[3,0] = 16.57643100960703
$arrEL[3,1] = 36.8900342768016
$arrEL[3,2] = 30.44525673169848
$arrEL[3,3] = 14.91859337076415
$arrEL[3,4] = 24.09824372694339
$arrEL[3,5] = 7.740028832514838
$arrEL[3,6] = 7.884189507983128
$arrEL[3,7] = 12.64721152796058
$arrEL[4,0] = 16.57658493077944
$arrEL[4,1] = 36.89322665764721
$arrEL[4,2] = 30.45294484060039
$arrEL[4,3] = 14.92109324159683
$arrEL[4,4] = 24.10326989110912
$arrEL[4,5] = 7.739497113563336
$arrEL[4,6] = 7.874307349482036
$arrEL[4,7] = 12.64615277689265
$arrEL[5,0] = 16.57574679238932
$arrEL[5,1] = 36.87180982212465
$arrEL[5,2] = 30.40032413430235
$arrEL[5,3] = 14.90392505509983
$arrEL[5,4] = 24.0687489534139
$arrEL[5,5] = 7.743183176059404
$arrEL[5,6] = 7.942676037794767
$arrEL[5,7] = 12.65365416919486
$arrEL[6,0] = 16.5777322325118
$arrEL[6,1] = 36.79567250208689
$arrEL[6,2] = 30.18655285793723
$arrEL[6,3] = 14.83270527693544
$arrEL[6,4] = 23.92546919645644
$arrEL[6,5] = 7.759330903313674
$arrEL[6,6] = 8.238259295291119
$arrEL[6,7] = 12.69049796927661
$arrEL[7,0] = 16.60190968374576
$arrEL[7,1] = 36.71249021209641
$arrEL[7,2] = 29.83400485215362
$arrEL[7,3] = 14.70931470026655
$arrEL[7,4] = 23.67697740760588
$arrEL[7,5] = 7.790543755231506
$arrEL[7,6] = 8.791016267746853
$arrEL[7,7] = 12.77587197975445
$arrEL[8,0] = 16.63189086000646
$arrEL[8,1] = 36.69172253362326
$arrEL[8,2] = 29.61599552631385
$arrEL[8,3] = 14.62853631073958
$arrEL[8,4] = 23.51415304056278
$arrEL[8,5] = 7.813175796279527
$arrEL[8,6] = 9.17685705497005
$arrEL[8,7] = 12.84628311025464
$arrEL[9,0] = 16.64814609368664
$arrEL[9,1] = 36.69102748308737
$arrEL[9,2] = 29.5258014848724
$arrEL[9,3] = 14.59392393480565
$arrEL[9,4] = 23.44435540317641
$arrEL[9,5] = 7.823401360490909
$arrEL[9,6] = 9.347273724104896
$arrEL[9,7] = 12.87991122042253
$arrEL[10,0] = 16.65467416678949
$arrEL[10,1] = 36.69202069109924
$arrEL[10,2] = 29.49294524120617
$arrEL[10,3] = 14.58112340305869
$arrEL[10,4] = 23.41853846487696
$arrEL[10,5] = 7.827262979982439
$arrEL[10,6] = 9.411020020426115
$arrEL[10,7] = 12.89286841922904
$arrEL[11,0] = 16.65325172125045
$arrEL[11,1] = 36.6917509502963
$arrEL[11,2] = 29.49996356800995
$arrEL[11,3] = 14.5838666054538
$arrEL[11,4] = 23.42407130632755
$arrEL[11,5] = 7.826431790239415
$arrEL[11,6] = 9.397326894483509
$arrEL[11,7] = 12.89006804916444
$arrEL[12,0] = 16.64867571754703
$arrEL[12,1] = 36.69108402999267
$arrEL[12,2] = 29.52307232222792
$arrEL[12,3] = 14.59286468873602
$arrEL[12,4] = 23.44221912649861
$arrEL[12,5] = 7.823719273365888
$arrEL[12,6] = 9.352534223153528
$arrEL[12,7] = 12.88097278826853
$arrEL[13,0] = 16.64592118799779
$arrEL[13,1] = 36.6908390649766
$arrEL[13,2] = 29.5373963844323
$arrEL[13,3] = 14.59841616664283
$arrEL[13,4] = 23.45341512896192
$arrEL[13,5] = 7.822056384707257
$arrEL[13,6] = 9.324993474548027
$arrEL[13,7] = 12.87543050877538
$arrEL[14,0] = 16.63088088582769
$arrEL[14,1] = 36.69194390987711
$arrEL[14,2] = 29.62207121228139
$arrEL[14,3] = 14.6308412151205
$arrEL[14,4] = 23.51880040880594
$arrEL[14,5] = 7.812506080338263
$arrEL[14,6] = 9.165612581818786
$arrEL[14,7] = 12.8441170104034
$arrEL[15,0] = 16.62232178868074
$arrEL[15,1] = 36.69486209759141
$arrEL[15,2] = 29.67632136286678
$arrEL[15,3] = 14.65127919833139
$arrEL[15,4] = 23.56000602255304
$arrEL[15,5] = 7.806629006875459
$arrEL[15,6] = 9.066491022734967
$arrEL[15,7] = 12.82531166620565
$arrEL[16,0] = 16.61764535284708
$arrEL[16,1] = 36.69736454984037
$arrEL[16,2] = 29.70836938809162
$arrEL[16,3] = 14.66323548542559
$arrEL[16,4] = 23.5841085216594
$arrEL[16,5] = 7.80324212603788
$arrEL[16,6] = 9.009001309453007
$arrEL[16,7] = 12.81464613950354
$arrEL[17,0] = 16.6161044420711
$arrEL[17,1] = 36.69835340970389
$arrEL[17,2] = 29.71936517425668
$arrEL[17,3] = 14.66731819782357
$arrEL[17,4] = 23.5923382909861
$arrEL[17,5] = 7.802094278632842
$arrEL[17,6] = 8.98945594376862
$arrEL[17,7] = 12.81106107964113
$arrEL[18,0] = 16.62320743013722
$arrEL[18,1] = 36.69446618480548
$arrEL[18,2] = 29.67045885969941
$arrEL[18,3] = 14.64908275114494
$arrEL[18,4] = 23.55557799972074
$arrEL[18,5] = 7.807255312457387
$arrEL[18,6] = 9.077092512108663
$arrEL[18,7] = 12.82729796727687
$arrEL[19,0] = 16.65000972022838
$arrEL[19,1] = 36.69124584350104
$arrEL[19,2] = 29.51624943002191
$arrEL[19,3] = 14.5902134208818
$arrEL[19,4] = 23.43687201331047
$arrEL[19,5] = 7.824516296742459
$arrEL[19,6] = 9.365712665455804
$arrEL[19,7] = 12.88363829111169
$arrEL[20,0] = 16.66969626982647
$arrEL[20,1] = 36.696462871035
$arrEL[20,2] = 29.42303527227782
$arrEL[20,3] = 14.55352478219081
$arrEL[20,4] = 23.36286898975331
$arrEL[20,5] = 7.835735428283587
$arrEL[20,6] = 9.549735628380567
$arrEL[20,7] = 12.92175621964839
$arrEL[21,0] = 16.65899195448725
$arrEL[21,1] = 36.69300944179624
$arrEL[21,2] = 29.47209039496147
$arrEL[21,3] = 14.57294293910717
$arrEL[21,4] = 23.40203850846245
$arrEL[21,5] = 7.829753410379688
$arrEL[21,6] = 9.45195678442167
$arrEL[21,7] = 12.90129566483732
$arrEL[22,0] = 16.62280627038614
$arrEL[22,1] = 36.69464260781764
$arrEL[22,2] = 29.67310662253157
$arrEL[22,3] = 14.65007512252389
$arrEL[22,4] = 23.55757862283251
$arrEL[22,5] = 7.806972184809075
$arrEL[22,6] = 9.072301142705239
$arrEL[22,7] = 12.82639950554151
$arrEL[23,0] = 16.59321163949836
$arrEL[23,1] = 36.72790749554631
$arrEL[23,2] = 29.9222039185675
$arrEL[23,3] = 14.74095765502182
$arrEL[23,4] = 23.74073001894731
$arrEL[23,5] = 7.782150772467683
$arrEL[23,6] = 8.644781010995926
$arrEL[23,7] = 12.75140224766681
$ws.Range("E2:L25").Value = $arrEL

$arrO = New-Object 'object[,]' 24,1
$arrO[0,0] = 22.69730888491643
$arrO[1,0] = 22.79766379066763
$arrO[2,0] = 22.86349462916699
$arrO[3,0] = 22.89138048889458
$arrO[4,0] = 22.89607489154111
$arrO[5,0] = 22.86386641872022
$arrO[6,0] = 22.73103693644569
$arrO[7,0] = 22.50398364252141
$arrO[8,0] = 22.3575447501298
$arrO[9,0] = 22.29535150215919
$arrO[10,0] = 22.2724367674926
$arrO[11,0] = 22.27734354810185
$arrO[12,0] = 22.29345353352048
$arrO[13,0] = 22.30340427242902
$arrO[14,0] = 22.36169826346199
$arrO[15,0] = 22.39859283111428
$arrO[16,0] = 22.42022984595661
$arrO[17,0] = 22.42762724033156
$arrO[18,0] = 22.39462225879095
$arrO[19,0] = 22.28870436138545
$arrO[20,0] = 22.22319117249166
$arrO[21,0] = 22.25781711415307
$arrO[22,0] = 22.39641602812394
$arrO[23,0] = 22.56182931086461
$ws.Range("O2:O25").Value = $arrO
